$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: PyTorch AutoGrad article title/link update
$ws.Range("D4").Value = "PyTorch의 자동미분(AutoGrad)을 활용하여 정형 데이터 예측기(regression model)구현"
$ws.Range("E4").Value = "https://teddylee777.github.io/pytorch/pytorch-boston-regression"

# Row 9: pabii blog title/link update
$ws.Range("D9").Value = "한국이 1류 국가로 도약할 길"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/why-korea-2nd-tier/#utm_source=rss&utm_medium=rss&utm_campaign=why-korea-2nd-tier"

# Row 51: bskyvision title/link update
$ws.Range("D51").Value = "웹 퍼블리셔와 프론트엔드 개발자의 차이는 무엇일까?"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EC%9B%B9-%ED%8D%BC%EB%B8%94%EB%A6%AC%EC%85%94%EC%99%80-%ED%94%84%EB%A1%A0%ED%8A%B8%EC%97%94%EB%93%9C-%EA%B0%9C%EB%B0%9C%EC%9E%90%EC%9D%98-%EC%B0%A8%EC%9D%B4%EB%8A%94-%EB%AC%B4%EC%97%87%EC%9D%BC%EA%B9%8C"
